$d = $word.ActiveDocument

# --- 1. Fix "Tracking " + "gantt" (two runs with spell-check proofErr
#        markers around the second run) into a single clean run
#        "Tracking gantt" with the proofErr elements removed. ---
#
# Directly overwriting the first paragraph's Range.Text leaves the
# second run's <w:proofErr/> markers behind as stray elements, so
# instead we insert a brand-new paragraph (which naturally gets a
# single plain run) right before the second paragraph, give it the
# merged text, and then delete the original (first) paragraph
# entirely - proofErr and all.

$pGanttChart = $d.Paragraphs.Item(2)            # "Gantt chart" paragraph
$pGanttChart.Range.InsertParagraphBefore() | Out-Null

$pNewTracking = $d.Paragraphs.Item(2)           # the freshly inserted paragraph
$pNewTracking.Range.Text = "Tracking gantt"

$d.Paragraphs.Item(1).Range.Delete()            # drop the old "Tracking "/"gantt" paragraph

# --- 2. Add three new list items after "Pert chart":
#        "Control flow diagram", "Sequence diagram", "Mobile app interface"
#        The document ends with a trailing empty list paragraph; insert
#        the new paragraphs right before it so they land after "Pert chart". ---

$newItems = @("Control flow diagram", "Sequence diagram", "Mobile app interface")
foreach ($itemText in $newItems) {
    $lastIndex = $d.Paragraphs.Count
    $pTrailing = $d.Paragraphs.Item($lastIndex)
    $pTrailing.Range.InsertParagraphBefore() | Out-Null
    $d.Paragraphs.Item($lastIndex).Range.Text = $itemText
}
